$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 16777215
$ws.Range("A1").Font.Size = 11
